$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PHM")

# Row 4 (Inventory)
$ws.Range("B4").Value = 8007000000.0
$ws.Range("C4").Value = -81265000.0
$ws.Range("D4").Value = -17513000.0
$ws.Range("E4").Value = 291130000.0
$ws.Range("F4").Value = -189364000.0

# Row 14 (Accounts Payable)
$ws.Range("B14").Value = 405000000.0
$ws.Range("C14").Value = 274578000.0
$ws.Range("D14").Value = 12940000.0
$ws.Range("E14").Value = -58959000.0
$ws.Range("F14").Value = -26910000.0

# Row 20 (Long Term Tax Liability (Deferred))
$ws.Range("B20").Value = -21000000.0
$ws.Range("C20").Value = 48106000.0
$ws.Range("D20").Value = 39831000.0
$ws.Range("E20").Value = 29706000.0
$ws.Range("F20").Value = 19955000.0

# Row 32 (Net Debt)
$ws.Range("G32").Value = 1840157000.0

# Row 33 (Total Debt)
$ws.Range("G33").Value = 3091613000.0
